$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 463, shifting the existing rows 463-469 down to 465-471.
$ws.Rows.Item(463).Insert()
$ws.Rows.Item(463).Insert()

# --- New row 463 ---
$ws.Cells.Item(463, 1).Value = 9
$ws.Cells.Item(463, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(463, 3).Value = "Metropolitana"
$ws.Cells.Item(463, 4).Value = 44595
$ws.Cells.Item(463, 5).Value = 13
$ws.Cells.Item(463, 6).Value = "Fruta"
$ws.Cells.Item(463, 7).Value = 100109
$ws.Cells.Item(463, 8).Value = "Uva"
$ws.Cells.Item(463, 9).Value = 100109001
$ws.Cells.Item(463, 10).Value = "Uva"
$ws.Cells.Item(463, 11).Value = "Flame Seedless"
$ws.Cells.Item(463, 12).Value = "Primera"
$ws.Cells.Item(463, 13).Value = 400
$ws.Cells.Item(463, 14).Value = 10000
$ws.Cells.Item(463, 15).Value = 10000
$ws.Cells.Item(463, 16).Value = 10000
$ws.Cells.Item(463, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(463, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(463, 19).Value = 556
$ws.Cells.Item(463, 20).Value = 18

# --- New row 464 ---
$ws.Cells.Item(464, 1).Value = 9
$ws.Cells.Item(464, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(464, 3).Value = "Metropolitana"
$ws.Cells.Item(464, 4).Value = 44595
$ws.Cells.Item(464, 5).Value = 13
$ws.Cells.Item(464, 6).Value = "Fruta"
$ws.Cells.Item(464, 7).Value = 100109
$ws.Cells.Item(464, 8).Value = "Uva"
$ws.Cells.Item(464, 9).Value = 100109001
$ws.Cells.Item(464, 10).Value = "Uva"
$ws.Cells.Item(464, 11).Value = "Superior Seedless"
$ws.Cells.Item(464, 12).Value = "Primera"
$ws.Cells.Item(464, 13).Value = 660
$ws.Cells.Item(464, 14).Value = 10000
$ws.Cells.Item(464, 15).Value = 11000
$ws.Cells.Item(464, 16).Value = 10576
$ws.Cells.Item(464, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(464, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(464, 19).Value = 588
$ws.Cells.Item(464, 20).Value = 18
